# Update the test data sheet:
#  - A5 (demo4@gmail.com) -> demo4@yahoo.com
#  - A6 (jay043patil@gmail.com) -> demo5@cts.com
#  - A7 (empty) -> jay043patil@gmail.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "demo4@yahoo.com"
$ws.Range("A6").Value = "demo5@cts.com"
$ws.Range("A7").Value = "jay043patil@gmail.com"
